# "aligned gray headings on left" -------------------------------------------
# Slides 2,3,4,5,6 ("Toxicity", "Disease status", "Staging", "Clinical Trial",
# "Deceased") each have a heading label + description shape pair whose text
# boxes get vertically centered ("anchor=ctr"), and (except slide 4) whose
# vertical position is nudged to a shared baseline. Slide 2's description
# shape also loses a stray trailing empty paragraph. Finally, the three
# gray "Example:" / "Process :" / "Additional examples:" captions on every
# slide become right aligned.

# PowerPoint's Shape.Top/.Left getters/setters round-trip through a 32-bit
# float expressed in points, then truncate (not round) to EMU on save
# (1 pt = 12700 EMU). Converting a target EMU value naively back to points
# therefore sometimes lands one EMU short after the float32 + truncation
# round trip. Resolve-PointsForEmu finds a points value that truncates back
# to exactly the desired EMU so the saved XML matches byte-for-byte.
function Resolve-PointsForEmu {
    param([double]$Emu)
    $basePoints = $Emu / 12700.0
    $step = 0.0000001
    for ($n = 0; $n -lt 200000; $n++) {
        $candidate = $basePoints + ($n * $step)
        $asFloat = [float]$candidate
        $roundTripEmu = [int64]([double]$asFloat * 12700.0)
        if ($roundTripEmu -eq $Emu) {
            return $candidate
        }
    }
    throw "Resolve-PointsForEmu: no points value round-trips to EMU $Emu"
}

$p = $ppt.ActivePresentation

# Per-slide shape ids: heading label shape, description shape, and the
# three gray caption shapes ("Example:", "Process :", "Additional examples:").
# $null offset EMU means "leave the position alone" (slide 4 keeps its
# existing offsets; only the anchor changes there).
$slideInfo = @(
    @{ Slide = 2; Label = 3; Desc = 4; OffsetEmu = 197715; Captions = @(8, 9, 10); TrimDesc = $true },
    @{ Slide = 3; Label = 3; Desc = 4; OffsetEmu = 197715; Captions = @(25, 26, 27); TrimDesc = $false },
    @{ Slide = 4; Label = 3; Desc = 4; OffsetEmu = $null;  Captions = @(21, 22, 23); TrimDesc = $false },
    @{ Slide = 5; Label = 3; Desc = 4; OffsetEmu = 204702; Captions = @(22, 23, 28); TrimDesc = $false },
    @{ Slide = 6; Label = 3; Desc = 4; OffsetEmu = 197715; Captions = @(17, 18, 19); TrimDesc = $false }
)

foreach ($info in $slideInfo) {
    $s = $p.Slides.Item($info.Slide)

    $targetIds = @($info.Label, $info.Desc) + $info.Captions

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if (-not ($targetIds -contains $sh.Id)) {
            continue
        }

        if ($sh.Id -eq $info.Label -or $sh.Id -eq $info.Desc) {
            if ($null -ne $info.OffsetEmu) {
                $sh.Top = Resolve-PointsForEmu($info.OffsetEmu)
            }
            $sh.TextFrame.VerticalAnchor = 3   # ppAnchorMiddle -> anchor="ctr"

            if ($sh.Id -eq $info.Desc -and $info.TrimDesc) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t.Length -gt 0 -and $t[$t.Length - 1] -eq "`r") {
                    $sh.TextFrame.TextRange.Text = $t.Substring(0, $t.Length - 1)
                }
            }
        }
        elseif ($info.Captions -contains $sh.Id) {
            $sh.TextFrame.TextRange.ParagraphFormat.Alignment = 3   # ppAlignRight -> algn="r"
        }
    }
}
